$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "year" column header becomes "season_ending_year"
$ws.Range("A1").Value = "season_ending_year"

# New column K: "calendar_year" - a numeric copy of the year column.
# Copy the formatting of the last header cell (J1) onto the new header (K1)
# so it keeps the same bold/centered/bordered header style.
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("K1").Value = "calendar_year"

# Numeric calendar_year values, one per data row, mirroring column A's year.
$ws.Range("K2").Value = 1960
$ws.Range("K3").Value = 1988
$ws.Range("K4").Value = 1972
$ws.Range("K5").Value = 1975
$ws.Range("K6").Value = 1966
